$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B154").Value = 53925
$ws.Range("E154").Value = 79.37
$ws.Range("F154").Value = 1
$ws.Range("G154").Value = 66.44

$ws.Range("B155").Value = 64350
$ws.Range("E155").Value = 70.63
$ws.Range("F155").Value = 101
$ws.Range("G155").Value = 6710.44

$ws.Range("B156").Value = 57756
$ws.Range("F156").Value = -100
$ws.Range("G156").Value = -6644

$ws.Range("B271").Value = 48706
$ws.Range("E271").Value = 39.8
$ws.Range("F271").Value = -144
$ws.Range("G271").Value = -4795.2

$ws.Range("B272").Value = 64973
$ws.Range("E272").Value = 35.4
$ws.Range("F272").Value = 150
$ws.Range("G272").Value = 4995

$ws.Range("B309").Value = 61610
$ws.Range("E309").Value = 122.71
$ws.Range("F309").Value = -58
$ws.Range("G309").Value = -5957.18

$ws.Range("B310").Value = 63565
$ws.Range("E310").Value = 109.19
$ws.Range("F310").Value = 60
$ws.Range("G310").Value = 6162.6

$ws.Range("B338").Value = 63520
$ws.Range("E338").Value = 153.4
$ws.Range("F338").Value = 97
$ws.Range("G338").Value = 13995.16

$ws.Range("B339").Value = 55373
$ws.Range("E339").Value = 163.62
$ws.Range("F339").Value = -94
$ws.Range("G339").Value = -13562.32

$ws.Range("B364").Value = 63652
$ws.Range("E364").Value = 55.42
$ws.Range("F364").Value = 250
$ws.Range("G364").Value = 13032.5

$ws.Range("B365").Value = 57885
$ws.Range("E365").Value = 62.28
$ws.Range("F365").Value = 4
$ws.Range("G365").Value = 208.52

$ws.Range("B367").Value = 63563
$ws.Range("E367").Value = 119.04
$ws.Range("F367").Value = 15
$ws.Range("G367").Value = 1679.4

$ws.Range("B368").Value = 61605
$ws.Range("E368").Value = 133.78
$ws.Range("F368").Value = -13
$ws.Range("G368").Value = -1455.48

$ws.Range("B374").Value = 63560
$ws.Range("E374").Value = 134.87
$ws.Range("F374").Value = 104
$ws.Range("G374").Value = 13193.44

$ws.Range("B375").Value = 60325
$ws.Range("E375").Value = 151.57
$ws.Range("F375").Value = -102
$ws.Range("G375").Value = -12939.72

$ws.Range("B381").Value = 57817
$ws.Range("F381").Value = 3
$ws.Range("G381").Value = 239.43

$ws.Range("B382").Value = 62865
$ws.Range("F382").Value = 151
$ws.Range("G382").Value = 12051.31

$ws.Range("B392").Value = 62933
$ws.Range("F392").Value = 146
$ws.Range("G392").Value = 8632.98

$ws.Range("B393").Value = 57835
$ws.Range("F393").Value = 1
$ws.Range("G393").Value = 59.13

$ws.Range("B411").Value = 57856
$ws.Range("F411").Value = 2
$ws.Range("G411").Value = 342.66

$ws.Range("B412").Value = 63007
$ws.Range("F412").Value = 984
$ws.Range("G412").Value = 168588.72

$ws.Range("B423").Value = 63102
$ws.Range("C423").Value = "HUL-Vim Bar Multipack Fw 4X200G"
$ws.Range("F423").Value = 36
$ws.Range("G423").Value = 2140.92

$ws.Range("B424").Value = 53082
$ws.Range("C424").Value = "HUL-VIM BAR MULTIPACK FW 4X200G"
$ws.Range("F424").Value = 1
$ws.Range("G424").Value = 59.47

$ws.Range("B528").Value = 58047
$ws.Range("D528").Value = 105.54
$ws.Range("E528").Value = 126.1
$ws.Range("F528").Value = 54
$ws.Range("G528").Value = 5699.16

$ws.Range("B529").Value = 47097
$ws.Range("D529").Value = 112.28
$ws.Range("E529").Value = 134.16
$ws.Range("F529").Value = 15
$ws.Range("G529").Value = 1684.2

$ws.Range("B578").Value = 45695
$ws.Range("E578").Value = 23.58
$ws.Range("F578").Value = -36
$ws.Range("G578").Value = -710.28

$ws.Range("B579").Value = 64915
$ws.Range("E579").Value = 20.98
$ws.Range("F579").Value = 40
$ws.Range("G579").Value = 789.2

$ws.Range("B596").Value = 65067
$ws.Range("E596").Value = 15.65
$ws.Range("F596").Value = 338
$ws.Range("G596").Value = 4978.74

$ws.Range("B597").Value = 53595
$ws.Range("E597").Value = 17.61
$ws.Range("F597").Value = -335
$ws.Range("G597").Value = -4934.55

$ws.Range("B679").Value = 53319
$ws.Range("E679").Value = 310.64
$ws.Range("F679").Value = -6
$ws.Range("G679").Value = -1643.52

$ws.Range("B680").Value = 64810
$ws.Range("E680").Value = 291.22
$ws.Range("F680").Value = 7
$ws.Range("G680").Value = 1917.44

$ws.Range("B701").Value = 64833
$ws.Range("E701").Value = 34.9
$ws.Range("F701").Value = 99
$ws.Range("G701").Value = 3250.17

$ws.Range("B702").Value = 60025
$ws.Range("E702").Value = 37.22
$ws.Range("F702").Value = -98
$ws.Range("G702").Value = -3217.34

$ws.Range("B712").Value = 64830
$ws.Range("E712").Value = 34.9
$ws.Range("F712").Value = 117
$ws.Range("G712").Value = 3841.11

$ws.Range("B713").Value = 60022
$ws.Range("E713").Value = 37.22
$ws.Range("F713").Value = -113
$ws.Range("G713").Value = -3709.79

$ws.Range("B864").Value = 65079
$ws.Range("E864").Value = 43.44
$ws.Range("F864").Value = 21
$ws.Range("G864").Value = 858.27

$ws.Range("B865").Value = 54751
$ws.Range("E865").Value = 46.34
$ws.Range("F865").Value = -19
$ws.Range("G865").Value = -776.53

Write-Output "Applied row swaps"
